$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calendar2021")
$ws.Activate()

$ws.Range("F4").Value = "HW 2"
$ws.Range("G4").Value = "HW 1"

$ws.Range("F5").Value = "PA 2"
$ws.Range("G5").Value = "HW 2;Quiz 0"

$ws.Range("F6").Value = "HW 3"

$ws.Range("G7").Value = " "

$ws.Range("F9").Value = ""
$ws.Range("G9").Value = ""

$ws.Range("G11").Value = ""

$ws.Range("G13").Value = ""

$ws.Range("G15").Value = ""

$ws.Range("G17").Value = ""

[void]$ws.Range("F22").Select()
